# Tidsplan opdateret for 21-03-2017
# Append new time-registration rows (55-61) for 21-03-2017 to the
# "Tidsregistrering" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

# Seed formatting for the new rows by copying it from the last existing
# "day" block (rows 49-54), so the new cells reuse the same style indexes
# (date format in column A, time format in columns G/H) instead of minting
# new ones.
$ws.Range("A49").Copy($ws.Range("A55"))
$ws.Range("G49:H50").Copy($ws.Range("G55:H56"))
$ws.Range("G51:H51").Copy($ws.Range("G57:H57"))
$ws.Range("G51:H51").Copy($ws.Range("G58:H58"))
$ws.Range("G51:H51").Copy($ws.Range("G59:H59"))
$ws.Range("G51:H51").Copy($ws.Range("G60:H60"))

# Row 55 - new date, role "Reviewer"
$ws.Range("A55").Value = [DateTime]"2017-03-21"
$ws.Range("E55").Value = "Reviewer"
$ws.Range("F55").Value = "Lavet review over OC12 design"
$ws.Range("G55").Value = 0.3576388888888889
$ws.Range("H55").Value = 0.47916666666666669

# Row 56 - role "Software Architect"
$ws.Range("E56").Value = "Software Architect"
$ws.Range("F56").Value = "Lavet OC14 design"
$ws.Range("G56").Value = 0.35416666666666669
$ws.Range("H56").Value = 0.47916666666666669

# Row 57
$ws.Range("F57").Value = "Lavet Test Suite for OC15: beregnSigmaB"
$ws.Range("G57").Value = 0.50694444444444442
$ws.Range("H57").Value = 0.52083333333333337

# Row 58
$ws.Range("F58").Value = "Lavet liste over klasser som skal Refactors "
$ws.Range("G58").Value = 0.52083333333333337
$ws.Range("H58").Value = 0.55208333333333337

# Row 59 - role "Reviewer"
$ws.Range("E59").Value = "Reviewer"
$ws.Range("F59").Value = "Lavet review over OC15 og OC16 design"
$ws.Range("G59").Value = 0.57291666666666663
$ws.Range("H59").Value = 0.60416666666666663

# Row 60
$ws.Range("F60").Value = "Rette Test Suite for OC15: beregnSigmaB"
$ws.Range("G60").Value = 0.60416666666666663
$ws.Range("H60").Value = 0.63194444444444442

# Row 61 - total hours for this day
$ws.Range("I61").Value = 5.3

# Scroll/selection state to mirror the saved view in the workbook (best
# effort - the headless host keeps the in-memory window position but the
# topLeftCell scroll offset itself isn't part of the exposed object model).
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H61").Select()
